$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 (HRNet / MXNET, OpenPose repo) FPS cell: fill in the missing benchmark
# number. The column stores these figures as plain text (matching the other
# FPS cells, e.g. D3="6.68", D7="11.16"), so build the text value via TEXT()
# in a scratch cell and paste just the value back in - this avoids Excel's
# "number stored as text" quote-prefix reformat (which would mint a new
# cell style) and keeps D6's existing style (the yellow-fill "TODO" format).
$scratch = $ws.Range("Z100")
$scratch.Formula = "=TEXT(5.28,""0.00"")"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)  # xlPasteValues
$scratch.Delete()

# Move the active selection to D2 (matches the saved selection state).
$ws.Range("D2").Select()
